# Auto-generated: apply numeric corrections to Chocobo_Profits workbook (scheduled price-data refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 24800
$ws.Range("J57").Value = 30000
$ws.Range("L57").Value = 90000
$ws.Range("N57").Value = -90998
$ws.Range("H64").Value = 3228.5715
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 3320
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 3320
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -3816
$ws.Range("H67").Value = 3228.5715
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 3320
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 3320
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -5036
$ws.Range("H100").Value = 15386984
$ws.Range("I100").Value = 16669085
$ws.Range("J100").Value = 1780
$ws.Range("K100").Value = 16669085
$ws.Range("L100").Value = 1780
$ws.Range("M100").Value = -16668544
$ws.Range("N100").Value = -2862
$ws.Range("H132").Value = 156477.22
$ws.Range("I132").Value = 2378.55
$ws.Range("K132").Value = 7135.650000000001
$ws.Range("M132").Value = -4605.650000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1211.0714
$ws.Range("I2").Value = 1205.375
$ws.Range("J2").Value = 1218.6666
$ws.Range("K2").Value = 1205.375
$ws.Range("L2").Value = 1218.6666
$ws.Range("M2").Value = -1092.375
$ws.Range("N2").Value = -1444.6666
$ws.Range("H32").Value = 8610.23
$ws.Range("I32").Value = 7602.778
$ws.Range("K32").Value = 7602.778
$ws.Range("M32").Value = -7315.778
$ws.Range("H44").Value = 35253.332
$ws.Range("J44").Value = 35253.332
$ws.Range("L44").Value = 35253.332
$ws.Range("N44").Value = -36229.332
$ws.Range("H45").Value = 1549
$ws.Range("I45").Value = 1542.8572
$ws.Range("K45").Value = 1542.8572
$ws.Range("M45").Value = -1165.8572
$ws.Range("H97").Value = 644.3333
$ws.Range("I97").Value = 661.2381
$ws.Range("J97").Value = 526
$ws.Range("K97").Value = 661.2381
$ws.Range("L97").Value = 526
$ws.Range("M97").Value = -165.2381
$ws.Range("N97").Value = -1518
$ws.Range("H116").Value = 1211.0714
$ws.Range("I116").Value = 1205.375
$ws.Range("J116").Value = 1218.6666
$ws.Range("K116").Value = 1205.375
$ws.Range("L116").Value = 1218.6666
$ws.Range("M116").Value = 1088.625
$ws.Range("N116").Value = -5806.6666
$ws.Range("H122").Value = 2037.8334
$ws.Range("I122").Value = 1265.2307
$ws.Range("J122").Value = 2950.9092
$ws.Range("K122").Value = 3795.6921
$ws.Range("L122").Value = 8852.7276
$ws.Range("M122").Value = -1345.6921
$ws.Range("N122").Value = -13752.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1211.0714
$ws.Range("I3").Value = 1205.375
$ws.Range("J3").Value = 1218.6666
$ws.Range("K3").Value = 1205.375
$ws.Range("L3").Value = 1218.6666
$ws.Range("M3").Value = -1091.375
$ws.Range("N3").Value = -1446.6666
$ws.Range("H94").Value = 656.63635
$ws.Range("I94").Value = 700.3333
$ws.Range("K94").Value = 700.3333
$ws.Range("M94").Value = -249.3333
$ws.Range("H99").Value = 1885.875
$ws.Range("I99").Value = 1125.375
$ws.Range("K99").Value = 1125.375
$ws.Range("M99").Value = 372.625
$ws.Range("H134").Value = 1810.1538
$ws.Range("I134").Value = 1245.3334
$ws.Range("J134").Value = 4916.6665
$ws.Range("K134").Value = 3736.0002
$ws.Range("L134").Value = 14749.9995
$ws.Range("M134").Value = -1201.0002
$ws.Range("N134").Value = -19819.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11365830
$ws.Range("I31").Value = 1303.4615
$ws.Range("J31").Value = 27781258
$ws.Range("K31").Value = 1303.4615
$ws.Range("L31").Value = 27781258
$ws.Range("M31").Value = -1008.4615
$ws.Range("N31").Value = -27781848
$ws.Range("H34").Value = 11365830
$ws.Range("I34").Value = 1303.4615
$ws.Range("J34").Value = 27781258
$ws.Range("K34").Value = 1303.4615
$ws.Range("L34").Value = 27781258
$ws.Range("M34").Value = -1101.4615
$ws.Range("N34").Value = -27781662
$ws.Range("H52").Value = 106800
$ws.Range("J52").Value = 106800
$ws.Range("L52").Value = 106800
$ws.Range("N52").Value = -107388
$ws.Range("H94").Value = 1700.0714
$ws.Range("I94").Value = 1799.6666
$ws.Range("J94").Value = 1672.909
$ws.Range("K94").Value = 1799.6666
$ws.Range("L94").Value = 1672.909
$ws.Range("M94").Value = -1348.6666
$ws.Range("N94").Value = -2574.909
$ws.Range("H99").Value = 11769923
$ws.Range("I99").Value = 28574386
$ws.Range("J99").Value = 6799
$ws.Range("K99").Value = 28574386
$ws.Range("L99").Value = 6799
$ws.Range("M99").Value = -28572888
$ws.Range("N99").Value = -9795
$ws.Range("H126").Value = 11769923
$ws.Range("I126").Value = 28574386
$ws.Range("J126").Value = 6799
$ws.Range("K126").Value = 85723158
$ws.Range("L126").Value = 20397
$ws.Range("M126").Value = -85720688
$ws.Range("N126").Value = -25337
$ws.Range("H139").Value = 40000
$ws.Range("J139").Value = 40000
$ws.Range("L139").Value = 40000
$ws.Range("N139").Value = -50280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 8298.857
$ws.Range("I69").Value = 2028
$ws.Range("J69").Value = 16660
$ws.Range("K69").Value = 6084
$ws.Range("L69").Value = 49980
$ws.Range("M69").Value = -5273
$ws.Range("N69").Value = -51602
$ws.Range("H72").Value = 8298.857
$ws.Range("I72").Value = 2028
$ws.Range("J72").Value = 16660
$ws.Range("K72").Value = 18252
$ws.Range("L72").Value = 149940
$ws.Range("M72").Value = -14196
$ws.Range("N72").Value = -158052

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 61.2
$ws.Range("I2").Value = 56.88889
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 56.88889
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = 56.11111
$ws.Range("N2").Value = -326
$ws.Range("H126").Value = 1827.25
$ws.Range("I126").Value = 1831.1753
$ws.Range("K126").Value = 5493.525900000001
$ws.Range("M126").Value = -3023.525900000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 6174466.5
$ws.Range("I93").Value = 10101954
$ws.Range("J93").Value = 2700
$ws.Range("K93").Value = 10101954
$ws.Range("L93").Value = 2700
$ws.Range("M93").Value = -10100706
$ws.Range("N93").Value = -5196

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 168484160
$ws.Range("J96").Value = 300001
$ws.Range("L96").Value = 300001
$ws.Range("N96").Value = -302747
$ws.Range("H107").Value = 1503
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1503
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 4509
$ws.Range("N107").Value = -8349
$ws.Range("M107").ClearContents()
